$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Successors (Test Results)")

# Update Big-U (E column) values for each successor block
$ws.Range("E1").Value = 0.8834590401531538
$ws.Range("E9").Value = 0.7739278245618743
$ws.Range("E17").Value = 0.6716430855176108
$ws.Range("E25").Value = 0.5766048230203634
$ws.Range("E33").Value = 0.4888130370701321
$ws.Range("E41").Value = 0.3873984592055312
$ws.Range("E49").Value = 0.2911532482517948
$ws.Range("E57").Value = 0.2156318289522959
$ws.Range("E65").Value = 0.156363822814376
$ws.Range("E73").Value = 0.1100783566303633

# Update Prev Op (H column) text for blocks 6-10
$ws.Range("H41").Value = "alloy_transform (in=[6.0, 12.0] out=[6.0, 6.0, 6.0]) (bins=6)"
$ws.Range("H49").Value = "alloy_transform (in=[7.0, 14.0] out=[7.0, 7.0, 7.0]) (bins=7)"
$ws.Range("H57").Value = "alloy_transform (in=[8.0, 16.0] out=[8.0, 8.0, 8.0]) (bins=8)"
$ws.Range("H65").Value = "alloy_transform (in=[9.0, 18.0] out=[9.0, 9.0, 9.0]) (bins=9)"
$ws.Range("H73").Value = "alloy_transform (in=[10.0, 20.0] out=[10.0, 10.0, 10.0]) (bins=10)"

# Update country rows (re-ordered alphabetically with new consumption values)
$ws.Range("A3").Value = "Atlantis"
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 68
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 11
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0

$ws.Range("A4").Value = "Brobdingnag"
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 30
$ws.Range("D4").Value = 12
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0

$ws.Range("A5").Value = "Carpania"
$ws.Range("B5").Value = 25
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 8
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0

$ws.Range("A6").Value = "Dinotopia"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 20
$ws.Range("D6").Value = 30
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0

$ws.Range("A7").Value = "Erewhon"
$ws.Range("B7").Value = 30
$ws.Range("C7").Value = 21
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0

$ws.Range("A11").Value = "Atlantis"
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = 66
$ws.Range("D11").Value = 20
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 11
$ws.Range("H11").Value = 2
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0

$ws.Range("A12").Value = "Brobdingnag"
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 30
$ws.Range("D12").Value = 12
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 10
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0

$ws.Range("A13").Value = "Carpania"
$ws.Range("B13").Value = 25
$ws.Range("C13").Value = 10
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 8
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0

$ws.Range("A14").Value = "Dinotopia"
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = 20
$ws.Range("D14").Value = 30
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 5
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0

$ws.Range("A15").Value = "Erewhon"
$ws.Range("B15").Value = 30
$ws.Range("C15").Value = 21
$ws.Range("D15").Value = 10
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0

$ws.Range("A19").Value = "Atlantis"
$ws.Range("B19").Value = 10
$ws.Range("C19").Value = 64
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = 8
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 11
$ws.Range("H19").Value = 3
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0

$ws.Range("A20").Value = "Brobdingnag"
$ws.Range("B20").Value = 5
$ws.Range("C20").Value = 30
$ws.Range("D20").Value = 12
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0

$ws.Range("A21").Value = "Carpania"
$ws.Range("B21").Value = 25
$ws.Range("C21").Value = 10
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 8
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0

$ws.Range("A22").Value = "Dinotopia"
$ws.Range("B22").Value = 3
$ws.Range("C22").Value = 20
$ws.Range("D22").Value = 30
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0

$ws.Range("A23").Value = "Erewhon"
$ws.Range("B23").Value = 30
$ws.Range("C23").Value = 21
$ws.Range("D23").Value = 10
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0

$ws.Range("A27").Value = "Atlantis"
$ws.Range("B27").Value = 10
$ws.Range("C27").Value = 62
$ws.Range("D27").Value = 20
$ws.Range("E27").Value = 9
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = 4
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0

$ws.Range("A28").Value = "Brobdingnag"
$ws.Range("B28").Value = 5
$ws.Range("C28").Value = 30
$ws.Range("D28").Value = 12
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 10
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0

$ws.Range("A29").Value = "Carpania"
$ws.Range("B29").Value = 25
$ws.Range("C29").Value = 10
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 8
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0

$ws.Range("A30").Value = "Dinotopia"
$ws.Range("B30").Value = 3
$ws.Range("C30").Value = 20
$ws.Range("D30").Value = 30
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0

$ws.Range("A31").Value = "Erewhon"
$ws.Range("B31").Value = 30
$ws.Range("C31").Value = 21
$ws.Range("D31").Value = 10
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 4
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0

$ws.Range("A35").Value = "Atlantis"
$ws.Range("B35").Value = 10
$ws.Range("C35").Value = 60
$ws.Range("D35").Value = 20
$ws.Range("E35").Value = 10
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 11
$ws.Range("H35").Value = 5
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0

$ws.Range("A36").Value = "Brobdingnag"
$ws.Range("B36").Value = 5
$ws.Range("C36").Value = 30
$ws.Range("D36").Value = 12
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 10
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0

$ws.Range("A37").Value = "Carpania"
$ws.Range("B37").Value = 25
$ws.Range("C37").Value = 10
$ws.Range("D37").Value = 3
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 8
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0

$ws.Range("A38").Value = "Dinotopia"
$ws.Range("B38").Value = 3
$ws.Range("C38").Value = 20
$ws.Range("D38").Value = 30
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 5
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0

$ws.Range("A39").Value = "Erewhon"
$ws.Range("B39").Value = 30
$ws.Range("C39").Value = 21
$ws.Range("D39").Value = 10
$ws.Range("E39").Value = 0
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 4
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0

$ws.Range("A43").Value = "Atlantis"
$ws.Range("B43").Value = 10
$ws.Range("C43").Value = 58
$ws.Range("D43").Value = 20
$ws.Range("E43").Value = 11
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 11
$ws.Range("H43").Value = 6
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0

$ws.Range("A44").Value = "Brobdingnag"
$ws.Range("B44").Value = 5
$ws.Range("C44").Value = 30
$ws.Range("D44").Value = 12
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 10
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0

$ws.Range("A45").Value = "Carpania"
$ws.Range("B45").Value = 25
$ws.Range("C45").Value = 10
$ws.Range("D45").Value = 3
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 8
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0

$ws.Range("A46").Value = "Dinotopia"
$ws.Range("B46").Value = 3
$ws.Range("C46").Value = 20
$ws.Range("D46").Value = 30
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 5
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0

$ws.Range("A47").Value = "Erewhon"
$ws.Range("B47").Value = 30
$ws.Range("C47").Value = 21
$ws.Range("D47").Value = 10
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 4
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0

$ws.Range("A51").Value = "Atlantis"
$ws.Range("B51").Value = 10
$ws.Range("C51").Value = 56
$ws.Range("D51").Value = 20
$ws.Range("E51").Value = 12
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 11
$ws.Range("H51").Value = 7
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0

$ws.Range("A52").Value = "Brobdingnag"
$ws.Range("B52").Value = 5
$ws.Range("C52").Value = 30
$ws.Range("D52").Value = 12
$ws.Range("E52").Value = 0
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 10
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0

$ws.Range("A53").Value = "Carpania"
$ws.Range("B53").Value = 25
$ws.Range("C53").Value = 10
$ws.Range("D53").Value = 3
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 8
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0

$ws.Range("A54").Value = "Dinotopia"
$ws.Range("B54").Value = 3
$ws.Range("C54").Value = 20
$ws.Range("D54").Value = 30
$ws.Range("E54").Value = 0
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 5
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0

$ws.Range("A55").Value = "Erewhon"
$ws.Range("B55").Value = 30
$ws.Range("C55").Value = 21
$ws.Range("D55").Value = 10
$ws.Range("E55").Value = 0
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 4
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0

$ws.Range("A59").Value = "Atlantis"
$ws.Range("B59").Value = 10
$ws.Range("C59").Value = 54
$ws.Range("D59").Value = 20
$ws.Range("E59").Value = 13
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 11
$ws.Range("H59").Value = 8
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0

$ws.Range("A60").Value = "Brobdingnag"
$ws.Range("B60").Value = 5
$ws.Range("C60").Value = 30
$ws.Range("D60").Value = 12
$ws.Range("E60").Value = 0
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 10
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0

$ws.Range("A61").Value = "Carpania"
$ws.Range("B61").Value = 25
$ws.Range("C61").Value = 10
$ws.Range("D61").Value = 3
$ws.Range("E61").Value = 0
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 8
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0

$ws.Range("A62").Value = "Dinotopia"
$ws.Range("B62").Value = 3
$ws.Range("C62").Value = 20
$ws.Range("D62").Value = 30
$ws.Range("E62").Value = 0
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 5
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0

$ws.Range("A63").Value = "Erewhon"
$ws.Range("B63").Value = 30
$ws.Range("C63").Value = 21
$ws.Range("D63").Value = 10
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 4
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0

$ws.Range("A67").Value = "Atlantis"
$ws.Range("B67").Value = 10
$ws.Range("C67").Value = 52
$ws.Range("D67").Value = 20
$ws.Range("E67").Value = 14
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 11
$ws.Range("H67").Value = 9
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0

$ws.Range("A68").Value = "Brobdingnag"
$ws.Range("B68").Value = 5
$ws.Range("C68").Value = 30
$ws.Range("D68").Value = 12
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 10
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0

$ws.Range("A69").Value = "Carpania"
$ws.Range("B69").Value = 25
$ws.Range("C69").Value = 10
$ws.Range("D69").Value = 3
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 8
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0

$ws.Range("A70").Value = "Dinotopia"
$ws.Range("B70").Value = 3
$ws.Range("C70").Value = 20
$ws.Range("D70").Value = 30
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 5
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0

$ws.Range("A71").Value = "Erewhon"
$ws.Range("B71").Value = 30
$ws.Range("C71").Value = 21
$ws.Range("D71").Value = 10
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 4
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0

$ws.Range("A75").Value = "Atlantis"
$ws.Range("B75").Value = 10
$ws.Range("C75").Value = 50
$ws.Range("D75").Value = 20
$ws.Range("E75").Value = 15
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 11
$ws.Range("H75").Value = 10
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0

$ws.Range("A76").Value = "Brobdingnag"
$ws.Range("B76").Value = 5
$ws.Range("C76").Value = 30
$ws.Range("D76").Value = 12
$ws.Range("E76").Value = 0
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 10
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0

$ws.Range("A77").Value = "Carpania"
$ws.Range("B77").Value = 25
$ws.Range("C77").Value = 10
$ws.Range("D77").Value = 3
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 8
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0

$ws.Range("A78").Value = "Dinotopia"
$ws.Range("B78").Value = 3
$ws.Range("C78").Value = 20
$ws.Range("D78").Value = 30
$ws.Range("E78").Value = 0
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 5
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0

$ws.Range("A79").Value = "Erewhon"
$ws.Range("B79").Value = 30
$ws.Range("C79").Value = 21
$ws.Range("D79").Value = 10
$ws.Range("E79").Value = 0
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 4
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
